# Correction in SA algorithm and 746 logs
# Update the "Fitness" column (C) for rows 2-104 of Sheet1 to reflect the
# corrected values produced by the algorithm fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 2-36 -> 7900
$ws.Range("C2:C36").Value = 7900

# Rows 37-103 -> 7748
$ws.Range("C37:C103").Value = 7748

# Row 104 -> 7295
$ws.Range("C104").Value = 7295
